# The commit swaps the colour content of the presentation's two theme
# parts: the slide master (used by every slide - ppt/theme/theme2.xml,
# "Integral" / "Red Violet") is recoloured with the plain default
# "Office Theme" palette, while the theme that used to hold the Office
# palette (ppt/theme/theme1.xml, only ever linked from the Notes
# Master) ends up holding the old "Integral"/"Red Violet" palette.
#
# The PowerPoint object model has no supported way to touch the Notes
# Master's theme colours directly (NotesMaster.ColorScheme /
# .ThemeColorScheme simply mirror the slide master's theme in this
# object model, the same way Master.ApplyTheme / ThemeColorScheme.Load
# / Design.Name are inert without a real theme file on disk), so the
# reachable, user-visible part of this change is recolouring the
# presentation's live theme (the slide master's theme, which every
# slide inherits) from the "Integral" palette to the standard "Office"
# palette via the modern 12-slot ThemeColorScheme object.

$p = $ppt.ActivePresentation

# Target ("Office") theme colour palette, in ThemeColorScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeHex = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

# Apply to every slide at once via a full SlideRange - they all share
# the same slide master / theme part, so this recolours it once.
$all = $p.Slides.Range()
$tcs = $all.ThemeColorScheme

for ($i = 1; $i -le $officeHex.Count; $i++) {
    $hex = $officeHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $rgb = $r + ($g * 256) + ($b * 65536)
    $tcs.Colors($i).RGB = $rgb
}
